# Update countries & provincias Spain
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4 - Estados Unidos
$ws.Range("B4").Value = 1376317
$ws.Range("C4").Value = 8679
$ws.Range("D4").Value = 258345
$ws.Range("E4").Value = 1036815
$ws.Range("G4").Value = 370
$ws.Range("H4").Value = 81157

# Row 9 - Francia
$ws.Range("B9").Value = 177423
$ws.Range("C9").Value = 453
$ws.Range("D9").Value = 56724
$ws.Range("E9").Value = 94056
$ws.Range("F9").Value = 2712
$ws.Range("G9").Value = 263
$ws.Range("H9").Value = 26643

# Row 60 - Barein
$ws.Range("F60").Value = 4

# Row 112 - Paraguay
$ws.Range("B112").Value = 724
$ws.Range("C112").Value = 11
$ws.Range("D112").Value = 170
$ws.Range("E112").Value = 544

# Row 173 - Malaui
$ws.Range("B173").Value = 57
$ws.Range("C173").Value = 1
$ws.Range("E173").Value = 40
